# Adds the "Furniture Sale Showcase Website" project to the PROJECTS sheet,
# as the newest/top entry (row 2), pushing the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PROJECTS")

# Insert a new blank row at row 2 (shifts existing data + blank filler rows down by one).
$ws.Rows("2:2").Insert()

# Fill in the new project's data.
$ws.Cells.Item(2, 1).Value = "Furniture Sale Showcase Website"
$ws.Cells.Item(2, 2).Value = "Web"
$ws.Cells.Item(2, 3).Value = "sell_furniture.png"
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "June 2025"
$ws.Cells.Item(2, 5).Value = "sell_furniture_i.png"
$ws.Cells.Item(2, 6).Value = "web"
$ws.Cells.Item(2, 7).Value = "html,css,react,bootstrap,javascript"
$ws.Cells.Item(2, 8).Value = "LINK|https://home-items-lluis.onrender.com/"
$ws.Cells.Item(2, 9).Value = "I created a small, responsive website to showcase furniture and items we planned to sell before moving out of an apartment. <br/> The site includes an image carousel for each item, item details like title and price, and a navigation system to browse through the collection. <br/> For privacy reasons, the images are not public, and the GitHub repository will remain private."

# Leave selection/activation on the PROJECTS sheet (matches the commit making it the active tab).
$ws.Activate()
$ws.Range("I3").Select()
